# Kjell Hansen reviews the "Reflektion" section and tightens up the
# final sentence of the report, leaving his edit as a tracked change
# (insertion/deletion) rather than a silent rewrite.

$d = $word.ActiveDocument

# Make sure the edit is attributed to the reviewer and recorded as a
# tracked change (Word records w:ins/w:del with the current
# Application.UserName as w:author).
$word.Application.UserName = "Kjell Hansen"
$d.TrackRevisions = $true

# Replace the second half of the sentence
#   "... förvirande. Detta är en dålig vana och som jag skall försöka
#   undvika i kommande projekt."
# with the reviewer's rewrite, keeping the leading "." of the previous
# sentence and the trailing "." intact so only the sentence itself is
# marked as changed.
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute(
    " Detta är en dålig vana och som jag skall försöka undvika i kommande projekt",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Detta är något jag kommer att tänka på när jag jobbar på andra projekt framöver",
    2)

# Word drops a zero-length "_GoBack" bookmark at the location of the
# last edit once the document is saved; recreate that marker right
# before the trailing period that follows the rewritten sentence.
$after = $d.Content
$after.Find.ClearFormatting()
$after.Find.Execute("andra projekt framöver")
$goBackPoint = $d.Range($after.End, $after.End)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $goBackPoint)
